$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell values (revised AgTests/AgPosit figures)
$ws.Range("F526").Value2 = 8965
$ws.Range("F527").Value2 = 11713
$ws.Range("F536").Value2 = 8179
$ws.Range("F537").Value2 = 13711
$ws.Range("F539").Value2 = 10704
$ws.Range("F540").Value2 = 12517
$ws.Range("F541").Value2 = 16757
$ws.Range("F542").Value2 = 10414
$ws.Range("F543").Value2 = 4756
$ws.Range("F544").Value2 = 14331
$ws.Range("F547").Value2 = 14008
$ws.Range("F548").Value2 = 17264
$ws.Range("F551").Value2 = 17883
$ws.Range("F554").Value2 = 17973
$ws.Range("F555").Value2 = 21624
$ws.Range("F556").Value2 = 12212
$ws.Range("F558").Value2 = 24696
$ws.Range("F559").Value2 = 22493
$ws.Range("F561").Value2 = 24334
$ws.Range("F562").Value2 = 27076
$ws.Range("F563").Value2 = 14094
$ws.Range("F565").Value2 = 28988
$ws.Range("F566").Value2 = 25738
$ws.Range("G566").Value2 = 322
$ws.Range("F567").Value2 = 23477
$ws.Range("F568").Value2 = 23968
$ws.Range("F569").Value2 = 32387
$ws.Range("F571").Value2 = 14988
$ws.Range("F572").Value2 = 33339
$ws.Range("F574").Value2 = 23336
$ws.Range("F575").Value2 = 25962
$ws.Range("F576").Value2 = 28841
$ws.Range("F577").Value2 = 14719
$ws.Range("F578").Value2 = 15031
$ws.Range("F579").Value2 = 32405
$ws.Range("F580").Value2 = 28692
$ws.Range("F581").Value2 = 26826
$ws.Range("F582").Value2 = 25673
$ws.Range("G582").Value2 = 469
$ws.Range("F583").Value2 = 28647
$ws.Range("G583").Value2 = 475
$ws.Range("F584").Value2 = 12928
$ws.Range("G584").Value2 = 245
$ws.Range("F585").Value2 = 14539
$ws.Range("G585").Value2 = 381

# Append new daily rows 586-588
$ws.Range("A586").NumberFormat = "yyyy-mm-dd"
$ws.Range("A586").Value2 = 44480
$ws.Range("B586").Value2 = 427480
$ws.Range("C586").Value2 = 11872
$ws.Range("D586").Value2 = 1509
$ws.Range("E586").Value2 = 12752
$ws.Range("F586").Value2 = 32537
$ws.Range("G586").Value2 = 679
$ws.Range("A587").NumberFormat = "yyyy-mm-dd"
$ws.Range("A587").Value2 = 44481
$ws.Range("B587").Value2 = 429886
$ws.Range("C587").Value2 = 13927
$ws.Range("D587").Value2 = 2406
$ws.Range("E587").Value2 = 12772
$ws.Range("F587").Value2 = 27110
$ws.Range("G587").Value2 = 534
$ws.Range("A588").NumberFormat = "yyyy-mm-dd"
$ws.Range("A588").Value2 = 44482
$ws.Range("B588").Value2 = 431757
$ws.Range("C588").Value2 = 12531
$ws.Range("D588").Value2 = 1871
$ws.Range("E588").Value2 = 12791
$ws.Range("F588").Value2 = 18533
$ws.Range("G588").Value2 = 387

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()